$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 160
$ws.Range("A160").Value = "'true"
$ws.Range("A160").ClearFormats()
$ws.Range("B160").Value = "'false"
$ws.Range("B160").ClearFormats()
$ws.Range("C160").Value = "'true"
$ws.Range("C160").ClearFormats()
$ws.Range("D160").Value = "NILL"
$ws.Range("E160").Value = "NILL"
$ws.Range("F160").Value = "NILL"
$ws.Range("G160").Value = "Kajal"
$ws.Range("H160").Value = "Sharma "
$ws.Range("I160").Value = "kvskajalsharma10@gmail.com"
$ws.Range("J160").Value = "'8010844492"
$ws.Range("J160").ClearFormats()
$ws.Range("K160").Value = "Kajal@123"
$ws.Range("L160").Value = "WiFUGhCVpmNtX1kBruZ2th6l4"
$ws.Range("M160").Value = "'0"
$ws.Range("M160").ClearFormats()
$ws.Range("N160").Value = "NILL"
$ws.Range("O160").Value = "NILL"
$ws.Range("P160").Value = "NILL"
$ws.Range("Q160").Value = "NILL"
$ws.Range("R160").Value = "NILL"
$ws.Range("S160").Value = "NILL"
$ws.Range("T160").Value = "NILL"
$ws.Range("U160").Value = "NILL"
$ws.Range("V160").Value = "'"
$ws.Range("V160").ClearFormats()
$ws.Range("W160").Value = "NILL"
$ws.Range("X160").Value = "NILL"
$ws.Range("Y160").Value = "https://rekonnectfileupload.s3.ap-south-1.amazonaws.com/RekonnectKajal%20Sharma%20resume.pdf"
$ws.Range("Z160").Value = "'"
$ws.Range("Z160").ClearFormats()
$ws.Range("AA160").Value = "'true"
$ws.Range("AA160").ClearFormats()
$ws.Range("AB160").Value = "Business development,Bach office work, I can do all types work from home "
$ws.Range("AC160").Value = "NILL"
$ws.Range("AD160").Value = "'true"
$ws.Range("AD160").ClearFormats()
$ws.Range("AE160").Value = "7days"

# Row 161
$ws.Range("A161").Value = "'true"
$ws.Range("A161").ClearFormats()
$ws.Range("B161").Value = "'false"
$ws.Range("B161").ClearFormats()
$ws.Range("C161").Value = "'true"
$ws.Range("C161").ClearFormats()
$ws.Range("D161").Value = "NILL"
$ws.Range("E161").Value = "NILL"
$ws.Range("F161").Value = "NILL"
$ws.Range("G161").Value = "Sarita"
$ws.Range("H161").Value = "Sarita"
$ws.Range("I161").Value = "sarita7july@gmail.com"
$ws.Range("J161").Value = "'9518168616"
$ws.Range("J161").ClearFormats()
$ws.Range("K161").Value = "snutrition2021"
$ws.Range("L161").Value = "dzcqtJSSIUSCDZdxRwcMldeu6"
$ws.Range("M161").Value = "'0"
$ws.Range("M161").ClearFormats()
$ws.Range("N161").Value = "NILL"
$ws.Range("O161").Value = "NILL"
$ws.Range("P161").Value = "NILL"
$ws.Range("Q161").Value = "NILL"
$ws.Range("R161").Value = "NILL"
$ws.Range("S161").Value = "NILL"
$ws.Range("T161").Value = "NILL"
$ws.Range("U161").Value = "NILL"
$ws.Range("V161").Value = "'"
$ws.Range("V161").ClearFormats()
$ws.Range("W161").Value = "NILL"
$ws.Range("X161").Value = "NILL"
$ws.Range("Y161").Value = "https://rekonnectfileupload.s3.ap-south-1.amazonaws.com/RekonnectSarita%20Biodata.pdf"
$ws.Range("Z161").Value = "'"
$ws.Range("Z161").ClearFormats()
$ws.Range("AA161").Value = "'false"
$ws.Range("AA161").ClearFormats()
$ws.Range("AB161").Value = "As a web developer"
$ws.Range("AC161").Value = "NILL"
$ws.Range("AD161").Value = "'true"
$ws.Range("AD161").ClearFormats()
$ws.Range("AE161").Value = "7days"

# Row 162
$ws.Range("A162").Value = "'true"
$ws.Range("A162").ClearFormats()
$ws.Range("B162").Value = "'false"
$ws.Range("B162").ClearFormats()
$ws.Range("C162").Value = "'true"
$ws.Range("C162").ClearFormats()
$ws.Range("D162").Value = "NILL"
$ws.Range("E162").Value = "NILL"
$ws.Range("F162").Value = "NILL"
$ws.Range("G162").Value = "Kumar"
$ws.Range("H162").Value = "Gaurav"
$ws.Range("I162").Value = "kgaurav.developer@gmail.com"
$ws.Range("J162").Value = "WlYeGWdam259oMdAdgK7tsw6u"
$ws.Range("K162").Value = "'0"
$ws.Range("K162").ClearFormats()
$ws.Range("L162").Value = "NILL"
$ws.Range("M162").Value = "NILL"
$ws.Range("N162").Value = "NILL"
$ws.Range("O162").Value = "NILL"
$ws.Range("P162").Value = "NILL"
$ws.Range("Q162").Value = "NILL"
$ws.Range("R162").Value = "NILL"
$ws.Range("S162").Value = "NILL"
$ws.Range("T162").Value = "'"
$ws.Range("T162").ClearFormats()
$ws.Range("U162").Value = "NILL"
$ws.Range("V162").Value = "NILL"
$ws.Range("W162").Value = "https://rekonnectfileupload.s3.ap-south-1.amazonaws.com/RekonnectGaurav%27s%20Resume.pdf"
$ws.Range("X162").Value = "'"
$ws.Range("X162").ClearFormats()
$ws.Range("Y162").Value = "'false"
$ws.Range("Y162").ClearFormats()
$ws.Range("Z162").Value = "I WANT TO BECOME A SUCCESSFUL ANDROID APP DEVELOPER"
$ws.Range("AA162").Value = "NILL"
$ws.Range("AB162").Value = "'true"
$ws.Range("AB162").ClearFormats()
$ws.Range("AC162").Value = "30days"

# Row 163
$ws.Range("A163").Value = "'true"
$ws.Range("A163").ClearFormats()
$ws.Range("B163").Value = "'false"
$ws.Range("B163").ClearFormats()
$ws.Range("C163").Value = "'true"
$ws.Range("C163").ClearFormats()
$ws.Range("D163").Value = "NILL"
$ws.Range("E163").Value = "NILL"
$ws.Range("F163").Value = "NILL"
$ws.Range("G163").Value = "Preksha"
$ws.Range("H163").Value = "Sethia"
$ws.Range("I163").Value = "prekshasethia4@gmail.com"
$ws.Range("J163").Value = "nKE4wnfkD9TTxRplggQ5TF6mO"
$ws.Range("K163").Value = "'0"
$ws.Range("K163").ClearFormats()
$ws.Range("L163").Value = "NILL"
$ws.Range("M163").Value = "NILL"
$ws.Range("N163").Value = "NILL"
$ws.Range("O163").Value = "NILL"
$ws.Range("P163").Value = "NILL"
$ws.Range("Q163").Value = "NILL"
$ws.Range("R163").Value = "NILL"
$ws.Range("S163").Value = "NILL"
$ws.Range("T163").Value = "'"
$ws.Range("T163").ClearFormats()
$ws.Range("U163").Value = "NILL"
$ws.Range("V163").Value = "NILL"
$ws.Range("W163").Value = "https://rekonnectfileupload.s3.ap-south-1.amazonaws.com/Rekonnect1624557057862Resume_Preksha.docx"
$ws.Range("X163").Value = "'"
$ws.Range("X163").ClearFormats()
$ws.Range("Y163").Value = "'false"
$ws.Range("Y163").ClearFormats()
$ws.Range("Z163").Value = "I m passionate about learning new things"
$ws.Range("AA163").Value = "NILL"
$ws.Range("AB163").Value = "'true"
$ws.Range("AB163").ClearFormats()
$ws.Range("AC163").Value = "15days"
